# "Generate Report for Handback" - update the localization-status workbook
# after a handback: refresh the Overview/zh-cn/de-de status text, stamp the
# handback datetimes, record the newly-generated target/handback files (with
# their hyperlinks), and widen the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$mdFileName = "9a68bde4-82d4-405d-a47f-0a482d644cba.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0864ddb089011d5b438db46e8a63c80bcb415472/e2e/9a68bde4-82d4-405d-a47f-0a482d644cba.md"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn row 2: Latest Target File / Latest Handback File / DateTime ---
$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$ws2.Range("J2").Value = "9a68bde4-82d4-405d-a47f-0a482d644cba.11f902fc3455ef1b669eb591a2a91b59a20edc7b.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-10-14 08:25:04"

# --- de-de row 2: Latest Target File / Latest Handback File / DateTime ---
$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$ws3.Range("J2").Value = "9a68bde4-82d4-405d-a47f-0a482d644cba.11f902fc3455ef1b669eb591a2a91b59a20edc7b.de-de.xlf"
$ws3.Range("K2").Value = "2016-10-14 08:25:21"

# --- Widen columns that now display the longer "Handed back..." status /
#     full-width file names (same pixel width the report generator used for
#     its other 40-char columns). ---
$ws1.Range("E1").EntireColumn.ColumnWidth = 29.17
$ws1.Range("F1").EntireColumn.ColumnWidth = 29.17

$ws2.Range("C1").EntireColumn.ColumnWidth = 29.17
$ws2.Range("I1").EntireColumn.ColumnWidth = 39.17
$ws2.Range("J1").EntireColumn.ColumnWidth = 39.17

$ws3.Range("C1").EntireColumn.ColumnWidth = 29.17
$ws3.Range("I1").EntireColumn.ColumnWidth = 39.17
$ws3.Range("J1").EntireColumn.ColumnWidth = 39.17
